$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a pure numeric-looking string must be forced to Text format
# first, otherwise Excel auto-converts the typed value into a Number (losing the exact
# literal text, e.g. trailing zeros) -- matching how the source workbook stores these
# as literal inline strings.
$numericLooking = @('D5', 'D6', 'D11', 'D13', 'D14', 'D15', 'D16', 'D22', 'D25', 'D30', 'D31', 'D32', 'D33', 'D36', 'D38', 'D41', 'D46', 'D48', 'D49', 'D51')
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new values cell-by-cell (values taken from the target diff)
$ws.Range('D2').Value = '68.906.80'
$ws.Range('E2').Value = '  -4.22%  '
$ws.Range('D3').Value = '3.502.15'
$ws.Range('E3').Value = '  -5.33%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '579.40'
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('D6').Value = '174.38'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').Value = '3.494.27'
$ws.Range('E8').Value = '  -5.28%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  -7.29%  '
$ws.Range('D11').Value = '6.63'
$ws.Range('E11').Value = '  +7.83%  '
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('D13').Value = '47.19'
$ws.Range('E13').Value = '  -5.56%  '
$ws.Range('D14').Value = '0.0000277'
$ws.Range('E14').Value = '  -3.56%  '
$ws.Range('D15').Value = '672.94'
$ws.Range('E15').Value = '  -2.23%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = '8.86'
$ws.Range('E16').Value = '  -1.44%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '4.064.78'
$ws.Range('E17').Value = '  -5.28%  '
$ws.Range('D18').Value = '3.508.91'
$ws.Range('E18').Value = '  -5.69%  '
$ws.Range('D19').Value = '68.847.09'
$ws.Range('E19').Value = '  -4.56%  '
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('E21').Value = '  -4.34%  '
$ws.Range('D22').Value = '11.17'
$ws.Range('E22').Value = '  -4.26%  '
$ws.Range('E23').Value = '  -4.52%  '
$ws.Range('E24').Value = '  -8.56%  '
$ws.Range('D25').Value = '98.15'
$ws.Range('E25').Value = '  -5.51%  '
$ws.Range('E26').Value = '  -3.91%  '
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -7.16%  '
$ws.Range('D30').Value = '9.42'
$ws.Range('E30').Value = '  -8.84%  '
$ws.Range('D31').Value = '32.89'
$ws.Range('E31').Value = '  -7.64%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '8.72'
$ws.Range('E32').Value = '  -5.40%  '
$ws.Range('B33').Value = 'Stacks'
$ws.Range('C33').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D33').Value = '3.20'
$ws.Range('E33').Value = '  -8.43%  '
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('E35').Value = '  -6.17%  '
$ws.Range('D36').Value = '577.20'
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('E37').Value = '  -3.63%  '
$ws.Range('D38').Value = '3.59'
$ws.Range('E38').Value = '  -14.74%  '
$ws.Range('E39').Value = '  -4.06%  '
$ws.Range('E40').Value = '  -5.39%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('E42').Value = '  -4.99%  '
$ws.Range('E43').Value = '  -5.02%  '
$ws.Range('E44').Value = '  -3.33%  '
$ws.Range('D45').Value = '3.423.74'
$ws.Range('E45').Value = '  -8.73%  '
$ws.Range('D46').Value = '33.30'
$ws.Range('E46').Value = '  -6.73%  '
$ws.Range('D47').Value = '0.0₃0704'
$ws.Range('E47').Value = '  -9.38%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').Value = '2.87'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').Value = '2.60'
$ws.Range('E49').Value = '  -7.50%  '
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').Value = '132.47'
$ws.Range('E51').Value = '  -1.32%  '
